$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# ------------------------------------------------------------------
# 1. Turn the "Call to action" run (in the btn-cta paragraph) into a
#    field-code based HYPERLINK ("placeholder for EN instructions"),
#    matching Word's own representation when a hyperlink is inserted
#    over pre-existing text: fldChar begin/instrText/separate/end
#    wrapping a Hyperlink-styled run that still reads "Call to action".
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Call to action`r") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Call to action' paragraph"
}

$insertionPoint = $target.Range.Duplicate
$insertionPoint.End = $insertionPoint.Start

$fieldXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:pPr><w:pStyle w:val="btn-cta"/></w:pPr>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:instrText>HYPERLINK "http://www.canada.ca/"</w:instrText></w:r>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr></w:r>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>
<w:r><w:rPr><w:color w:val="auto"/><w:highlight w:val="green"/><w:u w:val="none"/></w:rPr><w:t>Call to action</w:t></w:r>
<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>
</w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($fieldXml) | Out-Null

# Re-apply the "Hyperlink" character style to the visible result run
# (direct rStyle references inside raw InsertXML fragments get
# dropped, so it is (re)applied through the object model instead).
$full = $target.Range
$displayRun = $d.Range($full.End - 1 - 15, $full.End - 1)
$displayRun.Style = "Hyperlink"

# ------------------------------------------------------------------
# 2. Word keeps a single "_GoBack" bookmark marking the last edit
#    location: drop the one that used to sit after "Danger button"
#    and re-create it at the new edit point (end of the CTA
#    paragraph). This naturally renumbers "_Buttons" from 3 to 4,
#    exactly as in the target revision.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$goBackRange = $target.Range.Duplicate
$goBackRange.Start = $goBackRange.End - 1
$goBackRange.End = $goBackRange.End - 1
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
